# David Silva section correction again
#
# Fixes the typo in the "Lessons Learned" bullet on the David Joao slide:
#   "The dude how read in binary"  ->  "The dude how reads binary"
#
# The sentence is split the same way PowerPoint splits an edited run: the
# untouched leading text ("The dude how ") stays in its own run and the
# retyped tail ("reads binary") becomes a new run.

$p = $ppt.ActivePresentation

$needle = "The dude how read in binary"
$oldFragment = "read in binary"
$newFragment = "reads binary"

$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -like "*$needle*") {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

if ($targetShape -ne $null) {
    $fullRange = $targetShape.TextFrame.TextRange
    $fullRange.Replace($oldFragment, $newFragment) | Out-Null
}
